$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReleaseTracker")

# Colour constants (match the existing legend: green/blue/yellow/orange/red)
$colGreen  = 5296274   # RGB(146,208,80)  -> "cards added to tracker"
$colYellow = 65535     # RGB(255,255,0)   -> "waiting for images"

# --- Row status colour updates (Surging Sparks + promo ultra-rare checks) ---
# Most of these rows move from orange ("not out") to yellow ("waiting for images"),
# except rows 58 and 62 which move to green ("cards added to tracker").
$ws.Range("A55:D55").Interior.Color = $colYellow
$ws.Range("A56:D56").Interior.Color = $colYellow
$ws.Range("A57:D57").Interior.Color = $colYellow
$ws.Range("A58:D58").Interior.Color = $colGreen
$ws.Range("A59:D59").Interior.Color = $colYellow
$ws.Range("A60:D60").Interior.Color = $colYellow
$ws.Range("A61:D61").Interior.Color = $colYellow
$ws.Range("A62:D62").Interior.Color = $colGreen
$ws.Range("A63:D63").Interior.Color = $colYellow
$ws.Range("A64:D64").Interior.Color = $colYellow
$ws.Range("A65:D65").Interior.Color = $colYellow
$ws.Range("A67:D67").Interior.Color = $colYellow
$ws.Range("A68:D68").Interior.Color = $colYellow

# --- Notes (column E) updates ---
$ws.Cells.Item(55, 5).Value = "ultra rares added"
$ws.Cells.Item(56, 5).Value = "ultra rares added"
$ws.Cells.Item(61, 5).Value = "ultra rares added"
$ws.Cells.Item(62, 5).ClearContents()
$ws.Cells.Item(63, 5).Value = "ultra rares added, 4 store promos - check holo"
$ws.Cells.Item(65, 5).Value = "all 3 are pixel cosmos"
$ws.Cells.Item(67, 5).Value = "ultra rares added"
$ws.Cells.Item(68, 5).Value = "ultra rares added"

# --- Release date fix for the Mc Donnalds 2024 promo set ---
$ws.Cells.Item(69, 1).Value2 = 45630
